# Update BSAC yearly financials worksheet with latest reported figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("D8").Value = 3025900
$ws.Range("E8").Value = 3141500
$ws.Range("F8").Value = 3066400
$ws.Range("G8").Value = 3273700
$ws.Range("H8").Value = 2750700
$ws.Range("I8").Value = 2779700
$ws.Range("J8").Value = 2600000

# Row 15
$ws.Range("D15").Value = -114400
$ws.Range("E15").Value = -96100
$ws.Range("F15").Value = -78800
$ws.Range("G15").Value = -64900
$ws.Range("H15").Value = -89800
$ws.Range("I15").Value = -82900
$ws.Range("J15").Value = -78600

# Row 17
$ws.Range("D17").Value = 1515500
$ws.Range("E17").Value = 1762500
$ws.Range("F17").Value = 1829400
$ws.Range("G17").Value = 1859300
$ws.Range("H17").Value = 1703000
$ws.Range("I17").Value = 1840300
$ws.Range("J17").Value = 1635500

# Row 18
$ws.Range("D18").Value = 1510400
$ws.Range("E18").Value = 1379000
$ws.Range("F18").Value = 1237000
$ws.Range("G18").Value = 1414400
$ws.Range("H18").Value = 1047700
$ws.Range("I18").Value = 939400
$ws.Range("J18").Value = 964600

# Row 20
$ws.Range("D20").Value = -450700
$ws.Range("E20").Value = -523700
$ws.Range("F20").Value = -461700
$ws.Range("G20").Value = -492800
$ws.Range("H20").Value = -256100
$ws.Range("I20").Value = -343300
$ws.Range("J20").Value = -253300

# Row 21
$ws.Range("D21").Value = 1174100
$ws.Range("E21").Value = 951400
$ws.Range("F21").Value = 854200
$ws.Range("G21").Value = 986500
$ws.Range("H21").Value = 881400
$ws.Range("I21").Value = 679000
$ws.Range("J21").Value = 789900

# Row 23
$ws.Range("D23").Value = 1059700
$ws.Range("E23").Value = 855300
$ws.Range("F23").Value = 775300
$ws.Range("G23").Value = 921600
$ws.Range("H23").Value = 791600
$ws.Range("I23").Value = 596100
$ws.Range("J23").Value = 711300

# Row 24
$ws.Range("D24").Value = 211100
$ws.Range("E24").Value = 157500
$ws.Range("F24").Value = 110700
$ws.Range("G24").Value = 75000
$ws.Range("H24").Value = 138900
$ws.Range("I24").Value = 65300
$ws.Range("J24").Value = 113500

# Row 26
$ws.Range("D26").Value = 848600
$ws.Range("E26").Value = 697800
$ws.Range("F26").Value = 664600
$ws.Range("G26").Value = 846600
$ws.Range("H26").Value = 652800
$ws.Range("I26").Value = 530800
$ws.Range("J26").Value = 597800

# Row 27
$ws.Range("D27").Value = 830300
$ws.Range("E27").Value = 694400
$ws.Range("F27").Value = 659900
$ws.Range("G27").Value = 837800
$ws.Range("H27").Value = 649600
$ws.Range("I27").Value = 524000
$ws.Range("J27").Value = 590500

# Row 32
$ws.Range("D32").Value = 450700
$ws.Range("E32").Value = 523700
$ws.Range("F32").Value = 461700
$ws.Range("G32").Value = 492800
$ws.Range("H32").Value = 256100
$ws.Range("I32").Value = 343300
$ws.Range("J32").Value = 253300

# Row 33
$ws.Range("D33").Value = 830300
$ws.Range("E33").Value = 694400
$ws.Range("F33").Value = 659900
$ws.Range("G33").Value = 837800
$ws.Range("H33").Value = 649600
$ws.Range("I33").Value = 524000
$ws.Range("J33").Value = 590500

# Row 35
$ws.Range("D35").Value = 830300
$ws.Range("E35").Value = 694400
$ws.Range("F35").Value = 659900
$ws.Range("G35").Value = 837800
$ws.Range("H35").Value = 649600
$ws.Range("I35").Value = 524000
$ws.Range("J35").Value = 590500

# Row 41
$ws.Range("D41").Value = 2772100
$ws.Range("E41").Value = 2712900
$ws.Range("F41").Value = 2281900
$ws.Range("G41").Value = 3163700
$ws.Range("H41").Value = 3382900
$ws.Range("I41").Value = 2735800
$ws.Range("J41").Value = 4642000

# Row 42
$ws.Range("D42").Value = 4590400
$ws.Range("E42").Value = 6036900
$ws.Range("F42").Value = 7028000
$ws.Range("G42").Value = 5149100
$ws.Range("H42").Value = 2645200
$ws.Range("I42").Value = 507600
$ws.Range("J42").Value = 621400

# Row 47
$ws.Range("D47").Value = 40500
$ws.Range("E47").Value = 35000
$ws.Range("F47").Value = 29900
$ws.Range("G47").Value = 26300
$ws.Range("H47").Value = 14200
$ws.Range("I47").Value = 11200
$ws.Range("J47").Value = 12800

# Row 48
$ws.Range("D48").Value = 356500
$ws.Range("E48").Value = 378300
$ws.Range("F48").Value = 353800
$ws.Range("G48").Value = 311000
$ws.Range("H48").Value = 264900
$ws.Range("I48").Value = 238500
$ws.Range("J48").Value = 225000

# Row 49
$ws.Range("D49").Value = 92900
$ws.Range("E49").Value = 85400
$ws.Range("F49").Value = 75200
$ws.Range("G49").Value = 60200
$ws.Range("H49").Value = 98100
$ws.Range("I49").Value = 128400
$ws.Range("J49").Value = 118700

# Row 52
$ws.Range("D52").Value = 568200
$ws.Range("E52").Value = 550300
$ws.Range("F52").Value = 490400
$ws.Range("G52").Value = 417600
$ws.Range("H52").Value = 341100
$ws.Range("I52").Value = 267100
$ws.Range("J52").Value = 200700

# Row 54
$ws.Range("D54").Value = 52632300
$ws.Range("E54").Value = 54399800
$ws.Range("F54").Value = 50941500
$ws.Range("G54").Value = 44885200
$ws.Range("H54").Value = 39714700
$ws.Range("I54").Value = 36397000
$ws.Range("J54").Value = 36263400

# Row 57
$ws.Range("D57").Value = 715500
$ws.Range("E57").Value = 424100
$ws.Range("F57").Value = 679400
$ws.Range("G57").Value = 413500
$ws.Range("H57").Value = 406300
$ws.Range("I57").Value = 418900
$ws.Range("J57").Value = 131500

# Row 59
$ws.Range("D59").Value = 12200
$ws.Range("E59").Value = 45900
$ws.Range("F59").Value = 28600
$ws.Range("G59").Value = 4100
$ws.Range("H59").Value = 75600
$ws.Range("J59").Value = 2200

# Row 61
$ws.Range("D61").Value = 13280000
$ws.Range("E61").Value = 13939700
$ws.Range("F61").Value = 11003200
$ws.Range("G61").Value = 10616100
$ws.Range("H61").Value = 10394100
$ws.Range("I61").Value = 10518500
$ws.Range("J61").Value = 11518500

# Row 62
$ws.Range("D62").Value = 491000
$ws.Range("E62").Value = 465500
$ws.Range("F62").Value = 489500
$ws.Range("G62").Value = 467800
$ws.Range("H62").Value = 384100
$ws.Range("I62").Value = 296000
$ws.Range("J62").Value = 283500

# Row 66
$ws.Range("D66").Value = 48125000
$ws.Range("E66").Value = 50182800
$ws.Range("F66").Value = 46921500
$ws.Range("G66").Value = 41048700
$ws.Range("H66").Value = 36296000
$ws.Range("I66").Value = 33217300
$ws.Range("J66").Value = 33234200

# Row 72
$ws.Range("D72").Value = 3200500
$ws.Range("E72").Value = 2897000
$ws.Range("F72").Value = 2707900
$ws.Range("G72").Value = 2488700
$ws.Range("H72").Value = 2117300
$ws.Range("I72").Value = 1875100
$ws.Range("J72").Value = 1714900

# Row 76
$ws.Range("D76").Value = 4507300
$ws.Range("E76").Value = 4217000
$ws.Range("F76").Value = 4020000
$ws.Range("G76").Value = 3836500
$ws.Range("H76").Value = 3418700
$ws.Range("I76").Value = 3179800
$ws.Range("J76").Value = 3029300

# Row 81
$ws.Range("D81").Value = 830300
$ws.Range("E81").Value = 694400
$ws.Range("F81").Value = 659900
$ws.Range("G81").Value = 837800
$ws.Range("H81").Value = 649600
$ws.Range("I81").Value = 524000
$ws.Range("J81").Value = 590500

# Row 83
$ws.Range("D83").Value = 114400
$ws.Range("E83").Value = 96100
$ws.Range("F83").Value = 78800
$ws.Range("G83").Value = 64900
$ws.Range("H83").Value = 89800
$ws.Range("I83").Value = 82900
$ws.Range("J83").Value = 78600

# Row 89
$ws.Range("D89").Value = -589800
$ws.Range("E89").Value = 1105900
$ws.Range("F89").Value = 1026300
$ws.Range("G89").Value = 428200
$ws.Range("H89").Value = 787100
$ws.Range("I89").Value = -1696300
$ws.Range("J89").Value = 2288600

# Row 91
$ws.Range("D91").Value = -86400
$ws.Range("E91").Value = -91700
$ws.Range("F91").Value = -95700
$ws.Range("G91").Value = -86900
$ws.Range("H91").Value = -60000
$ws.Range("I91").Value = -54000
$ws.Range("J91").Value = -39200

# Row 94
$ws.Range("D94").Value = -108000
$ws.Range("E94").Value = -132600
$ws.Range("F94").Value = -136500
$ws.Range("G94").Value = -136200
$ws.Range("H94").Value = 44100
$ws.Range("I94").Value = -106000
$ws.Range("J94").Value = -76600

# Row 96
$ws.Range("D96").Value = -486000
$ws.Range("E96").Value = -494900
$ws.Range("F96").Value = -485400
$ws.Range("G96").Value = -389800
$ws.Range("H96").Value = -342200
$ws.Range("I96").Value = -383700
$ws.Range("J96").Value = -420900

# Row 100
$ws.Range("D100").Value = -508300
$ws.Range("E100").Value = -518600
$ws.Range("F100").Value = -500700
$ws.Range("G100").Value = -402800
$ws.Range("H100").Value = -192500
$ws.Range("I100").Value = -389900
$ws.Range("J100").Value = -425400

# Row 101
$ws.Range("D101").Value = -46200
$ws.Range("E101").Value = -220900
$ws.Range("F101").Value = 299100
$ws.Range("G101").Value = 51300
$ws.Range("H101").Value = -30400
$ws.Range("I101").Value = -5400
$ws.Range("J101").Value = -104600

# Row 102
$ws.Range("D102").Value = -1252200
$ws.Range("E102").Value = 233800
$ws.Range("F102").Value = 688200
$ws.Range("G102").Value = -59500
$ws.Range("H102").Value = 608300
$ws.Range("I102").Value = -2197600
$ws.Range("J102").Value = 1682000

